# "Zhe's manual labeling added (currently just mode 1)"
#  - Update a few "relatedness" values on Sheet1 per manual evaluation.
#  - Add a new worksheet "Sheet2" right after "Sheet1" containing a copy of
#    the table (with the updated relatedness labels), and make it the
#    active sheet with E11 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add a new sheet right after Sheet1 and name it Sheet2.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Copy the whole table (values + formats, so text-typed "prefix" numbers
# stay text and the header keeps its bold/border style) from Sheet1.
$ws1.Range("A1:E11").Copy()
$ws2.Range("A1:E11").PasteSpecial(-4104)

# Re-apply the header style explicitly (belt & braces in case a values-only
# paste ever changes) so the header keeps Sheet1's bold/border formatting.
$ws1.Range("A1:E1").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)

# Manual evaluation corrections to the "relatedness" column on Sheet2.
$ws2.Range("D3").Value = "Related"
$ws2.Range("D6").Value = "Related"
$ws2.Range("D7").Value = "Related"
$ws2.Range("D8").Value = "Slightly Related"
$ws2.Range("D10").Value = "Strongly Related"

# Sheet1 keeps its original data, just the selection resets to the table.
$ws1.Range("A1:E11").Select()

# Sheet2 becomes the active/selected sheet with the cursor on E11.
$ws2.Activate()
$ws2.Range("E11").Select()
